$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" '29.002.16'
Set-TextValue "E2" '  -0.74%  '
Set-TextValue "E3" '  -0.57%  '
Set-TextValue "D4" '0.9994'
Set-TextValue "E4" '  +0.08%  '
Set-TextValue "E5" '  -0.30%  '
Set-TextValue "D6" '0.6250'
Set-TextValue "E6" '  -5.90%  '
Set-TextValue "E7" '  +0.11%  '
Set-TextValue "B8" 'OKB'
Set-TextValue "C8" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D8" '45.19'
Set-TextValue "E8" '  +1.94%  '
Set-TextValue "B9" 'Dogecoin'
Set-TextValue "C9" 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue "D9" '0.07607'
Set-TextValue "E9" '  +2.16%  '
Set-TextValue "B10" 'Cardano'
Set-TextValue "C10" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D10" '0.2911'
Set-TextValue "E10" '  -1.43%  '
Set-TextValue "B11" 'Solana'
Set-TextValue "C11" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue "D11" '22.69'
Set-TextValue "E11" '  -2.50%  '
Set-TextValue "B12" 'TRON'
Set-TextValue "C12" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D12" '0.07751'
Set-TextValue "E12" '  -0.17%  '
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.834.72'
Set-TextValue "E13" '  -0.83%  '
Set-TextValue "B14" 'Polkadot'
Set-TextValue "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D14" '4.954'
Set-TextValue "E14" '  -1.44%  '
Set-TextValue "B15" 'Polygon'
Set-TextValue "C15" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D15" '0.6621'
Set-TextValue "E15" '  -1.72%  '
Set-TextValue "B16" 'Litecoin'
Set-TextValue "C16" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D16" '82.43'
Set-TextValue "E16" '  -1.15%  '
Set-TextValue "B17" 'ShibaInu'
Set-TextValue "C17" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.000009505'
Set-TextValue "E17" '  +10.09%  '
Set-TextValue "B18" 'Uniswap'
Set-TextValue "C18" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D18" '5.982'
Set-TextValue "E18" '  -3.10%  '
Set-TextValue "B19" 'WrappedBTC'
Set-TextValue "C19" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D19" '29.011.95'
Set-TextValue "E19" '  -0.65%  '
Set-TextValue "B20" 'BitcoinCash'
Set-TextValue "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D20" '223.68'
Set-TextValue "E20" '  -2.15%  '
Set-TextValue "B21" 'Avalanche'
Set-TextValue "C21" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D21" '12.32'
Set-TextValue "E21" '  -1.85%  '
Set-TextValue "B22" 'Dai'
Set-TextValue "C22" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D22" '1.000'
Set-TextValue "E22" '  +0.04%  '
Set-TextValue "B23" 'Chainlink'
Set-TextValue "C23" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D23" '7.205'
Set-TextValue "E23" '  +0.27%  '
Set-TextValue "B24" 'BinanceUSD'
Set-TextValue "C24" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D24" '1.001'
Set-TextValue "E24" '  +0.12%  '
Set-TextValue "B25" 'Monero'
Set-TextValue "C25" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D25" '158.78'
Set-TextValue "E25" '  -0.53%  '
Set-TextValue "B26" 'Cosmos'
Set-TextValue "C26" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D26" '8.393'
Set-TextValue "E26" '  -2.79%  '
Set-TextValue "B27" 'Stellar'
Set-TextValue "C27" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D27" '0.1362'
Set-TextValue "E27" '  -3.69%  '
Set-TextValue "B28" 'EthereumClassic'
Set-TextValue "C28" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D28" '17.83'
Set-TextValue "E28" '  -1.50%  '
Set-TextValue "B29" 'PancakeSwap'
Set-TextValue "C29" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D29" '1.492'
Set-TextValue "E29" '  -0.94%  '
Set-TextValue "B30" 'Filecoin'
Set-TextValue "C30" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D30" '4.053'
Set-TextValue "E30" '  -2.03%  '
Set-TextValue "B31" 'InternetComputer(DFINITY)'
Set-TextValue "C31" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D31" '4.023'
Set-TextValue "E31" '  -0.91%  '
Set-TextValue "B32" 'Toncoin'
Set-TextValue "C32" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D32" '1.192'
Set-TextValue "E32" '  +0.04%  '
Set-TextValue "B33" 'Hedera'
Set-TextValue "C33" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D33" '0.05182'
Set-TextValue "E33" '  -2.76%  '
Set-TextValue "B34" 'LidoDAOToken'
Set-TextValue "C34" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D34" '1.844'
Set-TextValue "E34" '  -2.11%  '
Set-TextValue "B35" 'ImmutableX'
Set-TextValue "C35" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D35" '0.7395'
Set-TextValue "E35" '  -0.47%  '
Set-TextValue "B36" 'ARBITRUM'
Set-TextValue "C36" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D36" '1.144'
Set-TextValue "E36" '  -1.09%  '
Set-TextValue "B37" 'HuobiToken'
Set-TextValue "C37" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D37" '2.698'
Set-TextValue "E37" '  +1.65%  '
Set-TextValue "B38" 'Maker'
Set-TextValue "C38" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D38" '1.254.83'
Set-TextValue "E38" '  -4.79%  '
Set-TextValue "B39" 'MXToken'
Set-TextValue "C39" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D39" '2.759'
Set-TextValue "E39" '  +0.60%  '
Set-TextValue "B40" 'VeChain'
Set-TextValue "C40" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D40" '0.01789'
Set-TextValue "E40" '  -0.65%  '
Set-TextValue "B41" 'FraxShare'
Set-TextValue "C41" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D41" '6.215'
Set-TextValue "E41" '  -3.30%  '
Set-TextValue "B42" 'TrustWalletToken'
Set-TextValue "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" '0.8919'
Set-TextValue "E42" '  -2.85%  '
Set-TextValue "B43" 'PaxDollar'
Set-TextValue "C43" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D43" '1.001'
Set-TextValue "E43" '  +0.18%  '
Set-TextValue "B44" 'Quant'
Set-TextValue "C44" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D44" '101.64'
Set-TextValue "E44" '  -1.53%  '
Set-TextValue "B45" 'RocketPoolETH'
Set-TextValue "C45" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D45" '1.976.36'
Set-TextValue "E45" '  -0.53%  '
Set-TextValue "B46" 'Aave'
Set-TextValue "C46" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D46" '64.44'
Set-TextValue "E46" '  -2.74%  '
Set-TextValue "B47" 'Mantle'
Set-TextValue "C47" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D47" '0.5108'
Set-TextValue "E47" '  -0.53%  '
Set-TextValue "B48" 'BabyDogeCoin'
Set-TextValue "C48" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D48" '0.00000000120'
Set-TextValue "E48" '  -0.13%  '
Set-TextValue "B49" 'TheSandbox'
Set-TextValue "C49" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D49" '0.3975'
Set-TextValue "E49" '  -1.18%  '
Set-TextValue "B50" 'EnergySwap'
Set-TextValue "C50" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D50" '8.830'
Set-TextValue "E50" '  +0.24%  '
Set-TextValue "B51" 'Cronos'
Set-TextValue "C51" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D51" '0.05752'
Set-TextValue "E51" '  -1.72%  '
